# B6-PowerPoint.pptx commit replay
#
# 1) Three tables (slides 14, 15, 16) get their table style swapped from
#    "{F468AA15-C224-4AE0-B284-AD96AEA04125}" to
#    "{B2D7E0AC-39CF-4BEF-B9AE-8AFDB70F7FAF}".
# 2) The deck's theme colour scheme (the "Integral" / Red Violet palette)
#    is replaced by the default Office Theme colour palette.

$p = $ppt.ActivePresentation

# --- 1. Retarget the table style on the three affected slides ---------
$oldStyle = "{F468AA15-C224-4AE0-B284-AD96AEA04125}"
$newStyle = "{B2D7E0AC-39CF-4BEF-B9AE-8AFDB70F7FAF}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyle) {
                $tbl.ApplyStyle($newStyle)
            }
        }
    }
}

# --- 2. Swap the theme colour scheme back to the default Office colours
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72
